$d = $word.ActiveDocument

$replacements = @(
    @("81÷7=11, 4", "72÷9=8, 0"),
    @("59÷5=11, 4", "35÷9=3, 8"),
    @("43÷5=8, 3", "88÷9=9, 7"),
    @("34÷2=17, 0", "35÷4=8, 3"),
    @("77÷9=8, 5", "80÷5=16, 0"),
    @("58÷6=9, 4", "20÷4=5, 0"),
    @("63÷9=7, 0", "74÷6=12, 2"),
    @("99÷6=16, 3", "55÷2=27, 1"),
    @("95÷4=23, 3", "95÷5=19, 0"),
    @("72÷7=10, 2", "47÷9=5, 2"),
    @("62÷6=10, 2", "16÷3=5, 1"),
    @("18÷3=6, 0", "34÷4=8, 2"),
    @("98÷5=19, 3", "51÷8=6, 3"),
    @("27÷2=13, 1", "99÷5=19, 4"),
    @("43÷9=4, 7", "17÷6=2, 5"),
    @("51÷5=10, 1", "36÷7=5, 1"),
    @("92÷4=23, 0", "37÷8=4, 5"),
    @("54÷4=13, 2", "88÷5=17, 3"),
    @("11÷3=3, 2", "13÷2=6, 1"),
    @("82÷2=41, 0", "58÷2=29, 0"),
    @("59÷3=19, 2", "54÷8=6, 6"),
    @("22÷3=7, 1", "76÷9=8, 4"),
    @("55÷8=6, 7", "49÷3=16, 1"),
    @("85÷8=10, 5", "87÷4=21, 3"),
    @("71÷6=11, 5", "55÷3=18, 1")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
